$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 90
$ws.Cells.Item($row, 1).Value = "2024-10-31 00:00:00"
$ws.Cells.Item($row, 2).Value = 73700
$ws.Cells.Item($row, 3).Value = 10322.13
$ws.Cells.Item($row, 4).Value = 9134.629999999999
$ws.Cells.Item($row, 5).Value = 7.1186
